# Add federal tax ("Bundessteuer") data for tax year 2018 to the second
# worksheet, mirroring the structure of the existing 2012 data block, then
# leave the UI focused on the newly-entered range (matches the author's
# end-of-edit selection/active-sheet state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Jahr, Steuertyp, Tariftyp, Einkommen, Steuer, Steuerinkrement, Einkommensinkrement
$newRows = @(
  ,@(2018, 1, 1, 17800,  25.41,     0.77,  100)
  ,@(2018, 1, 1, 31700,  132.53,    0.88,  100)
  ,@(2018, 1, 1, 41500,  220.54,    2.64,  100)
  ,@(2018, 1, 1, 55300,  585.17,    2.97,  100)
  ,@(2018, 1, 1, 72600,  1101.94,   5.94,  100)
  ,@(2018, 1, 1, 78200,  1435.2,    6.6,   100)
  ,@(2018, 1, 1, 103700, 3120.4,    8.8,   100)
  ,@(2018, 1, 1, 134700, 5850.6,    11,    100)
  ,@(2018, 1, 1, 176100, 10406.8,   13.2,  100)
  ,@(2018, 1, 1, 755300, 86859.5,   11.5,  100)
  ,@(2018, 1, 2, 30800,  25,        1,     100)
  ,@(2018, 1, 2, 51000,  228,       2,     100)
  ,@(2018, 1, 2, 58500,  379,       3,     100)
  ,@(2018, 1, 2, 75400,  887,       4,     100)
  ,@(2018, 1, 2, 90400,  1488,      5,     100)
  ,@(2018, 1, 2, 103500, 2144,      6,     100)
  ,@(2018, 1, 2, 114800, 2823,      7,     100)
  ,@(2018, 1, 2, 124300, 3489,      8,     100)
  ,@(2018, 1, 2, 131800, 4090,      9,     100)
  ,@(2018, 1, 2, 137400, 4595,      10,    100)
  ,@(2018, 1, 2, 141300, 4986,      11,    100)
  ,@(2018, 1, 2, 143200, 5196,      12,    100)
  ,@(2018, 1, 2, 145100, 5425,      13,    100)
  ,@(2018, 1, 2, 895900, 103028.5,  11.5,  100)
)

$startRow = 26
for ($i = 0; $i -lt $newRows.Count; $i++) {
  $r = $startRow + $i
  $vals = $newRows[$i]
  for ($c = 0; $c -lt $vals.Count; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $vals[$c]
  }
}

# Switch focus to the Bundessteuer sheet and select the newly added A-column
# block, with the view scrolled so row 33 is at the top - matches the
# sheetView/selection state captured in the saved workbook.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 33
$win.ScrollColumn = 1
$ws.Range("A26:A49").Select()
